$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 161815.39
$ws.Range("J12").Value = 300373
$ws.Range("L12").Value = 300373
$ws.Range("N12").Value = -300713
$ws.Range("H21").Value = 23890.777
$ws.Range("I21").Value = 27504.25
$ws.Range("K21").Value = 27504.25
$ws.Range("M21").Value = -27036.25
$ws.Range("H23").Value = 23890.777
$ws.Range("I23").Value = 27504.25
$ws.Range("K23").Value = 27504.25
$ws.Range("M23").Value = -27270.25
$ws.Range("H62").Value = 2637.4666
$ws.Range("I62").Value = 2137.625
$ws.Range("J62").Value = 3208.7144
$ws.Range("K62").Value = 2137.625
$ws.Range("L62").Value = 3208.7144
$ws.Range("M62").Value = -1513.625
$ws.Range("N62").Value = -4456.7144
$ws.Range("H65").Value = 2637.4666
$ws.Range("I65").Value = 2137.625
$ws.Range("J65").Value = 3208.7144
$ws.Range("K65").Value = 10688.125
$ws.Range("L65").Value = 16043.572
$ws.Range("M65").Value = -7568.125
$ws.Range("N65").Value = -22283.572
$ws.Range("H129").Value = 6411292.5
$ws.Range("I129").Value = 62500348
$ws.Range("J129").Value = 1114.6571
$ws.Range("K129").Value = 187501044
$ws.Range("L129").Value = 3343.9713
$ws.Range("M129").Value = -187496044
$ws.Range("N129").Value = -13343.9713
$ws.Range("H132").Value = 20007000
$ws.Range("I132").Value = 25006252
$ws.Range("J132").Value = 9998
$ws.Range("K132").Value = 75018756
$ws.Range("L132").Value = 29994
$ws.Range("M132").Value = -75016226
$ws.Range("N132").Value = -35054
$ws.Range("H133").Value = 23653.334
$ws.Range("J133").Value = 23653.334
$ws.Range("L133").Value = 23653.334
$ws.Range("N133").Value = -33773.334
$ws.Range("H134").Value = 23846.875
$ws.Range("J134").Value = 23846.875
$ws.Range("L134").Value = 23846.875
$ws.Range("N134").Value = -33986.875
$ws.Range("H135").Value = 532.4194
$ws.Range("I135").Value = 416.83334
$ws.Range("J135").Value = 4000
$ws.Range("K135").Value = 3751.50006
$ws.Range("L135").Value = 36000
$ws.Range("M135").Value = -1216.50006
$ws.Range("N135").Value = -41070
$ws.Range("H136").Value = 29642.857
$ws.Range("J136").Value = 29642.857
$ws.Range("L136").Value = 29642.857
$ws.Range("N136").Value = -39842.857
$ws.Range("H137").Value = 4444.0967
$ws.Range("I137").Value = 4504
$ws.Range("J137").Value = 4238.7144
$ws.Range("K137").Value = 13512
$ws.Range("L137").Value = 12716.1432
$ws.Range("M137").Value = -10962
$ws.Range("N137").Value = -17816.1432
$ws.Range("H138").Value = 3712.942
$ws.Range("I138").Value = 1908.5349
$ws.Range("J138").Value = 5517.3486
$ws.Range("K138").Value = 5725.6047
$ws.Range("L138").Value = 16552.0458
$ws.Range("M138").Value = -585.6046999999999
$ws.Range("N138").Value = -26832.0458
$ws.Range("H139").Value = 23812.5
$ws.Range("J139").Value = 23812.5
$ws.Range("L139").Value = 23812.5
$ws.Range("N139").Value = -34092.5
$ws.Range("H140").Value = 28905.715
$ws.Range("J140").Value = 28905.715
$ws.Range("L140").Value = 28905.715
$ws.Range("N140").Value = -39265.715
$ws.Range("H141").Value = 445838.22
$ws.Range("I141").Value = 1018.5
$ws.Range("J141").Value = 2225117.2
$ws.Range("K141").Value = 3055.5
$ws.Range("L141").Value = 6675351.600000001
$ws.Range("M141").Value = 2124.5
$ws.Range("N141").Value = -6685711.600000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14212.235
$ws.Range("I32").Value = 12630.754
$ws.Range("J32").Value = 22407.182
$ws.Range("K32").Value = 12630.754
$ws.Range("L32").Value = 22407.182
$ws.Range("M32").Value = -12343.754
$ws.Range("N32").Value = -22981.182

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 25859.762
$ws.Range("I86").Value = 1219.375
$ws.Range("K86").Value = 1219.375
$ws.Range("M86").Value = -96.375
$ws.Range("H89").Value = 25859.762
$ws.Range("I89").Value = 1219.375
$ws.Range("K89").Value = 6096.875
$ws.Range("M89").Value = -480.875

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3782.0417
$ws.Range("I62").Value = 2456.1177
$ws.Range("J62").Value = 7002.143
$ws.Range("K62").Value = 2456.1177
$ws.Range("L62").Value = 7002.143
$ws.Range("M62").Value = -1832.1177
$ws.Range("N62").Value = -8250.143
$ws.Range("H65").Value = 3782.0417
$ws.Range("I65").Value = 2456.1177
$ws.Range("J65").Value = 7002.143
$ws.Range("K65").Value = 12280.5885
$ws.Range("L65").Value = 35010.715
$ws.Range("M65").Value = -9160.588499999998
$ws.Range("N65").Value = -41250.715
$ws.Range("H105").Value = 2625.04
$ws.Range("I105").Value = 2524.2222
$ws.Range("K105").Value = 2524.2222
$ws.Range("M105").Value = -777.2222000000002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 248
$ws.Range("I17").Value = 248
$ws.Range("K17").Value = 744
$ws.Range("M17").Value = -575
$ws.Range("H56").Value = 5531.4614
$ws.Range("I56").Value = 5531.4614
$ws.Range("K56").Value = 5531.4614
$ws.Range("M56").Value = -5001.4614
$ws.Range("H63").Value = 9740
$ws.Range("I63").Value = 2051.4285
$ws.Range("J63").Value = 17428.572
$ws.Range("K63").Value = 6154.2855
$ws.Range("L63").Value = 52285.716
$ws.Range("M63").Value = -5405.2855
$ws.Range("N63").Value = -53783.716
$ws.Range("H66").Value = 9740
$ws.Range("I66").Value = 2051.4285
$ws.Range("J66").Value = 17428.572
$ws.Range("K66").Value = 18462.8565
$ws.Range("L66").Value = 156857.148
$ws.Range("M66").Value = -14718.8565
$ws.Range("N66").Value = -164345.148
$ws.Range("H70").Value = 3876.5
$ws.Range("I70").Value = 3012
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 9036
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -8721
$ws.Range("N70").Value = -12630
$ws.Range("H73").Value = 3876.5
$ws.Range("I73").Value = 3012
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 9036
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -7944
$ws.Range("N73").Value = -14184
$ws.Range("H103").Value = 2683.0344
$ws.Range("I103").Value = 1562.5
$ws.Range("J103").Value = 2862.32
$ws.Range("K103").Value = 4687.5
$ws.Range("L103").Value = 8586.960000000001
$ws.Range("M103").Value = -3808.5
$ws.Range("N103").Value = -10344.96
$ws.Range("H113").Value = 1165.5625
$ws.Range("J113").Value = 1165.5625
$ws.Range("L113").Value = 3496.6875
$ws.Range("N113").Value = -7836.6875
$ws.Range("H114").Value = 675.7826
$ws.Range("I114").Value = 198.5
$ws.Range("J114").Value = 1196.4546
$ws.Range("K114").Value = 595.5
$ws.Range("L114").Value = 3589.3638
$ws.Range("M114").Value = 2658.5
$ws.Range("N114").Value = -10097.3638
$ws.Range("H117").Value = 1360.0416
$ws.Range("I117").Value = 436
$ws.Range("J117").Value = 2452.0908
$ws.Range("K117").Value = 1308
$ws.Range("L117").Value = 7356.2724
$ws.Range("M117").Value = 2134
$ws.Range("N117").Value = -14240.2724

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 20803.8
$ws.Range("I36").Value = 666.6667
$ws.Range("J36").Value = 51009.5
$ws.Range("K36").Value = 666.6667
$ws.Range("L36").Value = 51009.5
$ws.Range("M36").Value = -181.6667
$ws.Range("N36").Value = -51979.5
$ws.Range("H57").Value = 18530.5
$ws.Range("H132").Value = 5381.926
$ws.Range("I132").Value = 7684
$ws.Range("K132").Value = 23052
$ws.Range("M132").Value = -20522

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H138").Value = 59800
$ws.Range("J138").Value = 59800
$ws.Range("L138").Value = 59800
$ws.Range("N138").Value = -70080
$ws.Range("H139").Value = 40544.453
$ws.Range("J139").Value = 40544.453
$ws.Range("L139").Value = 40544.453
$ws.Range("N139").Value = -50824.453
$ws.Range("H140").Value = 29616.125
$ws.Range("J140").Value = 29616.125
$ws.Range("L140").Value = 29616.125
$ws.Range("N140").Value = -39976.125

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H86").Value = 20125
$ws.Range("J86").Value = 20125
$ws.Range("L86").Value = 20125
$ws.Range("N86").Value = -22371
$ws.Range("H89").Value = 20125
$ws.Range("J89").Value = 20125
$ws.Range("L89").Value = 100625
$ws.Range("N89").Value = -111857
$ws.Range("H132").Value = 21676.68
$ws.Range("I132").Value = 5980.7
$ws.Range("J132").Value = 84460.60000000001
$ws.Range("K132").Value = 17942.1
$ws.Range("L132").Value = 253381.8
$ws.Range("M132").Value = -15412.1
$ws.Range("N132").Value = -258441.8
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140
$ws.Range("H136").Value = 2815.4194
$ws.Range("I136").Value = 2262.6365
$ws.Range("J136").Value = 4166.6665
$ws.Range("K136").Value = 6787.9095
$ws.Range("L136").Value = 12499.9995
$ws.Range("M136").Value = -4237.9095
$ws.Range("N136").Value = -17599.9995
$ws.Range("H137").Value = 40328.57
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 40328.57
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 40328.57
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -50528.57
$ws.Range("H138").Value = 30912.5
$ws.Range("J138").Value = 30912.5
$ws.Range("L138").Value = 30912.5
$ws.Range("N138").Value = -41192.5
$ws.Range("H139").Value = 44933.332
$ws.Range("J139").Value = 44933.332
$ws.Range("L139").Value = 44933.332
$ws.Range("N139").Value = -55213.332

Write-Output "Applied all cell updates"